$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.567.24"
$ws.Range("E2").Value = "  +4.67%  "
$ws.Range("D3").Value = "3.367.42"
$ws.Range("E3").Value = "  +9.63%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +8.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.387"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "3.365.27"
$ws.Range("E10").Value = "  +9.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "98.156.78"
$ws.Range("E13").Value = "  +4.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.73%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000248"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.997.10"
$ws.Range("E16").Value = "  +9.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.54%  "
$ws.Range("D18").Value = "3.369.22"
$ws.Range("E18").Value = "  +9.78%  "
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000211"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("E25").Value = "  +3.71%  "
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "3.543.87"
$ws.Range("E28").Value = "  +9.47%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("E31").Value = "  +4.88%  "
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "525.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.450"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.21%  "
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.785"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +16.82%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  +5.98%  "
$ws.Range("E49").Value = "  +6.80%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.72%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.10%  "
